$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TaxCode")
$ws.Activate()

# Header row
$ws.Range("A1").Value = "Tax Code"
$ws.Range("B1").Value = "Tax Percent"
$ws.Range("C1").Value = "Description"
$ws.Range("A1:C1").Font.Bold = $true

# Data rows (written in this order so new shared-string entries line up
# with how the workbook was originally authored)
$ws.Range("A3").Value = "INVALIDTC"
$ws.Range("B3").Value = "invalidpercentage"
$ws.Range("C2").Value = "Valid test data"
$ws.Range("C3").Value = "Invalid test data"
$ws.Range("C4").Value = "Update test data"
$ws.Range("A2").Value = "ZC"
$ws.Range("B2").Value = 12
$ws.Range("B4").Value = 13

# Column widths
$ws.Columns.Item(1).ColumnWidth = 14.67
$ws.Columns.Item(2).ColumnWidth = 15.17
$ws.Columns.Item(3).ColumnWidth = 13.33

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("A2").Select()
